$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Dillon Brooks", "SG,SF", "Houston Rockets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Draymond Green", "PF,C", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 6
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
